$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "Promedio" (average) column to both contribution tables.
# Table 1: header row 2, data rows 3-12   -> column R
# Table 2: header row 16, data rows 17-26 -> column R
# ---------------------------------------------------------------------------

# --- Table 1 header (row 2) --------------------------------------------------
$ws.Range("Q2").Copy()
$ws.Range("R2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("R2").Value = "Promedio"

# --- Table 1 data (rows 3-12) ------------------------------------------------
$t1First = $ws.Range("R3")
$t1First.Formula = "=AVERAGE(B3:Q3)"
$t1First.Borders.LineStyle = 1
$t1First.NumberFormat = "0.0"

for ($r = 4; $r -le 12; $r++) {
    $t1First.Copy()
    $dest = $ws.Range("R" + $r)
    $dest.PasteSpecial(-4122)   # xlPasteFormats
    $dest.Formula = "=AVERAGE(B" + $r + ":Q" + $r + ")"
}

# --- Table 2 header (row 16) --------------------------------------------------
$ws.Range("Q16").Copy()
$ws.Range("R16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R16").Value = "Promedio"

# --- Table 2 data (rows 17-26) ------------------------------------------------
$t1First.Copy()
$t2First = $ws.Range("R17")
$t2First.PasteSpecial(-4122)   # xlPasteFormats
$t2First.Formula = "=AVERAGE(B17:Q17)"

for ($r = 18; $r -le 26; $r++) {
    $t1First.Copy()
    $dest = $ws.Range("R" + $r)
    $dest.PasteSpecial(-4122)   # xlPasteFormats
    $dest.Formula = "=AVERAGE(B" + $r + ":Q" + $r + ")"
}

$excel.CutCopyMode = 0
